$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.824.34'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '3.215.49'
$ws.Range("E3").Value = '  -2.36%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.213.10'
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("E11").Value = '  -3.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("D15").Value = '3.744.98'
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").Value = '3.211.28'
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").Value = '63.754.70'
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.708'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("E23").Value = '  -2.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("E34").Value = '  -2.40%  '
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").Value = '0.0₃0731'
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0394'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '409.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("E43").Value = '  -3.84%  '
$ws.Range("D44").Value = '2.828.80'
$ws.Range("E44").Value = '  -8.83%  '
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  -0.32%  '
